$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column data: G2 = "Gold", H2 = 10
$ws.Range("G2").Value = "Gold"
$ws.Range("H2").Value = 10

# Update existing Amount values in column D
$ws.Range("D3").Value = 6
$ws.Range("D4").Value = 11
$ws.Range("D5").Value = 29

# Update the selected cell to match the target view state
$ws.Range("K2").Select()
